$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Helper pattern used throughout: to append text right before a
# paragraph's trailing paragraph mark, shrink the range by one
# character (to exclude the mark), collapse to its end, then
# InsertAfter the new text. Re-fetching Paragraphs(idx).Range each time
# keeps the range in sync with the document after each mutation.
# ---------------------------------------------------------------------

# --- Paragraph 4: "Some new creatures and a tweak on the Vanilla ones." ---
# Add a new sentence (preceded by a line break) about nocturnal predators.
$r = $d.Paragraphs(4).Range
$r.End = $r.End - 1
$r.Collapse(0)
$r.InsertAfter([char]11 + "First off, note that all predators are nocturnal animals. They sleep by day light and hunt by night.")

# --- New paragraph (becomes paragraph 5) about herbivores/omnivores ---
$d.Paragraphs(4).Range.InsertParagraphAfter()

$r = $d.Paragraphs(5).Range
$r.InsertAfter("Herbivores and omnivores on the other hand are diurnal animals. They sleep by night.")
$r = $d.Paragraphs(5).Range
$r.End = $r.End - 1
$r.Collapse(0)
$r.InsertAfter([char]11 + "Understand that sleeping or not carnivores will attack if approached within 5 blocks. ")

# --- New paragraph (becomes paragraph 6): mod changes behavior of creatures ---
$d.Paragraphs(5).Range.InsertParagraphAfter()

$r = $d.Paragraphs(6).Range
$r.InsertAfter("The mod changes the behavior of all the creatures in the game.")
$r = $d.Paragraphs(6).Range
$r.End = $r.End - 1
$r.Collapse(0)
$r.InsertAfter([char]11)

# --- New paragraph (becomes paragraph 7): The Wolf ---
$d.Paragraphs(6).Range.InsertParagraphAfter()

$r = $d.Paragraphs(7).Range
$r.InsertAfter("The Wolf:")

$r = $d.Paragraphs(7).Range
$r.End = $r.End - 1
$r.Collapse(0)
$r.InsertAfter([char]11 + "Wolves come in 9 different colors, but one of the things that the mod adds is pelts, and now, when you kill a wolf you can skin it for the same pelt color. ")

$r = $d.Paragraphs(7).Range
$r.End = $r.End - 1
$r.Collapse(0)
$r.InsertAfter("Also")

$r = $d.Paragraphs(7).Range
$r.End = $r.End - 1
$r.Collapse(0)
$r.InsertAfter(" is possible to get a head trophy of the same color.")

$r = $d.Paragraphs(7).Range
$r.End = $r.End - 1
$r.Collapse(0)
$r.InsertAfter([char]11)

# --- New paragraph (becomes paragraph 8): The Black Panther ---
$d.Paragraphs(7).Range.InsertParagraphAfter()

$r = $d.Paragraphs(8).Range
$r.InsertAfter("The Black Panther:")

$r = $d.Paragraphs(8).Range
$r.End = $r.End - 1
$r.Collapse(0)
$r.InsertAfter([char]11 + "It has similar behavior as the wolf, ")

$r = $d.Paragraphs(8).Range
$r.End = $r.End - 1
$r.Collapse(0)
$r.InsertAfter("and also")

$r = $d.Paragraphs(8).Range
$r.End = $r.End - 1
$r.Collapse(0)
$r.InsertAfter(" can give you a nice pelt and or a head trophy.")

$r = $d.Paragraphs(8).Range
$r.End = $r.End - 1
$r.Collapse(0)
$r.InsertAfter([char]11)

# --- "Hyenas:" is inserted directly into the existing bookmark
# paragraph (paragraph 9), right before the bookmark, instead of
# becoming its own paragraph. ---
$r = $d.Paragraphs(9).Range
$r.Collapse(1)
$r.InsertBefore("Hyenas:")

# --- Append three empty paragraphs after the bookmark paragraph, before sectPr ---
$r = $d.Paragraphs(9).Range
$r.InsertParagraphAfter()
$r = $d.Paragraphs(9).Range
$r.InsertParagraphAfter()
$r = $d.Paragraphs(9).Range
$r.InsertParagraphAfter()
